$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = 0.7525028612459097
$ws.Range("J5").Value = 0.4834040075754473
$ws.Range("K5").Value = 0.05804491740751398
$ws.Range("L5").Value = 2.279527057949995
